# Reshaped text box for team member names
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)   # "Subtitle 2" placeholder (subTitle, idx=1)

# Resize / reposition the subtitle placeholder (values converted from EMU to points, 1 pt = 12700 EMU)
$shape.Left   = 1295400 / 12700   # 102.0 pt
$shape.Top    = 3886200 / 12700   # 306.0 pt
$shape.Width  = 6629400 / 12700   # 522.0 pt
$shape.Height = 1752600 / 12700   # 138.0 pt

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Merge the ", " + "Lydia, " runs into a single ", Lydia, " run
$full = $tr.Text
$target1 = ", Lydia, "
$idx1 = $full.IndexOf($target1)
$chars1 = $tr.Characters($idx1 + 1, $target1.Length)
$chars1.Text = $target1

# Merge the " Fein" + ", Kevin McKenna, Bob Stark, Matthew " runs into one run
$full2 = $tr.Text
$target2 = " Fein, Kevin McKenna, Bob Stark, Matthew "
$idx2 = $full2.IndexOf($target2)
$chars2 = $tr.Characters($idx2 + 1, $target2.Length)
$chars2.Text = $target2
